# Practice: Solve Leetcode Problem#155 Min Stack
# Append two new rows (27 and 28) to the tracking sheet for this problem,
# mirroring the existing "Anna" / "Stephan" row-pair pattern used
# throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 27 (Anna's entry) ----
$ws.Cells.Item(27, 1).Value2 = "LeetCode"
$ws.Cells.Item(27, 1).HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(27, 2).Value2 = "Anna"
$ws.Cells.Item(27, 2).HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(27, 3).Value2 = "Easy"
$ws.Cells.Item(27, 4).Value2 = "155. Min Stack"

# ---- Row 28 (Stephan's entry) ----
$ws.Cells.Item(28, 1).Value2 = "LeetCode"
$ws.Cells.Item(28, 1).HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(28, 2).Value2 = "Stephan"
$ws.Cells.Item(28, 2).HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(28, 4).Value2 = "155. Min Stack"
$ws.Cells.Item(28, 4).WrapText = $true
$ws.Cells.Item(28, 5).Value2 = "2020/12/16"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 6).Value2 = "Stack"
$ws.Cells.Item(28, 7).Value2 = "Completed"

# ---- Update the saved selection/active cell, as Excel records on save ----
$ws.Range("F34").Select() | Out-Null
